# Update weekly excess mortality
# Mirrors the OOXML diff for "Berekening oversterfte CBS.xlsx":
#  - revise a handful of previously-entered "Waargenomen" (observed) weekly
#    death counts in column G (weeks 16, 33, 34, 36-42)
#  - append a brand-new week 43 row (row 35: F=week, G=observed, H=expected,
#    I=G-H oversterfte)
#  - move the view: scroll so row 8 is at the top and select L19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- revised observed values (column G) ------------------------------
$ws.Range("G8").Value  = 4304   # week 16: 4303 -> 4304
$ws.Range("G25").Value = 3207   # week 33: 3206 -> 3207
$ws.Range("G26").Value = 2846   # week 34: 2845 -> 2846
$ws.Range("G28").Value = 2682   # week 36: 2679 -> 2682
$ws.Range("G29").Value = 2734   # week 37: 2733 -> 2734
$ws.Range("G30").Value = 2712   # week 38: 2710 -> 2712
$ws.Range("G31").Value = 2883   # week 39: 2881 -> 2883
$ws.Range("G32").Value = 2993   # week 40: 2988 -> 2993
$ws.Range("G33").Value = 3005   # week 41: 2996 -> 3005
$ws.Range("G34").Value = 3197   # week 42: 3224 -> 3197

# I-column holds "=G-H" (oversterfte); it is a shared formula already
# covering this range so it recalculates automatically for rows 8 & 25-34.

# --- new row 35: week 43 ----------------------------------------------
$ws.Range("F35").Value = 43
$ws.Range("G35").Value = 3452
$ws.Range("H35").Value = 2862
$ws.Range("I35").Formula = "=G35-H35"

# G37/H37/I37 are SUM formulas and recalc automatically with the above.

# --- view/selection state ----------------------------------------------
# Scroll the window so row 8 is the top-left visible row (was row 2), then
# move the selection from I38 to L19.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L19").Select()
